# Edit script: add Csets/Sets (column B) to Commodities and Processes sheets,
# and reorder the Commodities rows (4-46) to match the target layout.

$wb = $excel.ActiveWorkbook

# ---- Commodities sheet ----
$wsComm = $wb.Worksheets.Item("Commodities")

$commB = @('NRG', 'ENV', 'ENV', 'NRG', 'ENV', 'MAT', 'MAT', 'DEM', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'NRG', 'ENV', 'NRG', 'NRG', 'NRG', 'ENV', 'MAT', 'MAT', 'NRG', 'NRG', 'MAT', 'NRG', 'NRG', 'ENV', 'MAT', 'NRG', 'NRG', 'MAT', 'MAT', 'NRG', 'ENV', 'NRG')
$commC = @('sec_biogas', 'emi_CO2_f_ind', 'emi_CO2_f_x2x_neg_reusable', 'pri_uran', 'emi_N2O_f_ind', 'iip_steel_oxygen', 'iip_steel_iron_pellets', 'exo_steel', 'pri_deuterium', 'pri_coal', 'iip_steel_blafu_slag', 'pri_crude_oil', 'pri_hydro_energy', 'sec_methane', 'sec_heating_oil', 'sec_hydrogen', 'sec_heat_low', 'iip_heat_proc', 'sec_elec_ind', 'pri_biomass', 'pri_waste', 'CO2_f_pow', 'emi_CO2_f_x2x_neg_stored]', 'iip_coke', 'pri_solar_radiation', 'sec_heavy_fuel_oil', 'emi_CH4_f_ind', 'iip_steel_raw_iron', 'iip_steel_crudesteel', 'pri_geoth_heat', 'pri_wind_energy_on', 'iip_steel_scrap', 'sec_elec', 'sec_natural_gas_syn', '[emi_CO2_f_x2x_neg_reusable', 'iip_steel_iron_ore', 'pri_natural_gas', 'pri_wind_energy_off', 'iip_steel_sinter', 'iip_steel_sponge_iron', 'sec_heat_high', 'emi_CO2_f_x2x', 'sec_H2')

for ($i = 0; $i -lt $commB.Length; $i++) {
    $row = 4 + $i
    $wsComm.Cells.Item($row, 2).Value = $commB[$i]
    $wsComm.Cells.Item($row, 3).Value = $commC[$i]
}

# ---- Processes sheet ----
$wsProc = $wb.Worksheets.Item("Processes")

$procB = @('PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'DEM', 'DEM', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'CHP', 'PRE', 'CHP', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE', 'PRE')

for ($i = 0; $i -lt $procB.Length; $i++) {
    $row = 4 + $i
    $wsProc.Cells.Item($row, 2).Value = $procB[$i]
}

